# Apply the edit described by the diff:
# 1. Rename the worksheet/sheet tab from "GossF-HW25.xpc" to "GossF"
# 2. Fix tiny floating point roundings in row 13 (D13, J13, K13, L13)
# 3. Append a new row 16 with results for the "HexGrid-60degTilt5degRes" case

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "GossF"

# 2. Tiny precision corrections on row 13
$ws.Range("D13").Value = 0.9943698965328298
$ws.Range("J13").Value = 0.9943698965328298
$ws.Range("K13").Value = 0.9942114746140635
$ws.Range("L13").Value = 0.9944075333387624

# 3. New row 16 (HexGrid-60degTilt5degRes, index 13 / 14th entry)
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9898537016774713
$ws.Range("D16").Value = 0.9491496470760641
$ws.Range("E16").Value = 1.018040081265707
$ws.Range("F16").Value = 0.9898537016774713
$ws.Range("G16").Value = 0.9478784471131123
$ws.Range("H16").Value = 1.086662413002611
$ws.Range("I16").Value = 1.010712039757009
$ws.Range("J16").Value = 0.9491496470760641
$ws.Range("K16").Value = 0.9835948641708856
$ws.Range("L16").Value = 0.9867242829241785
$ws.Range("M16").Value = 1.000382721648662
